$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 318, shifting existing rows 318:336 down to 319:337.
$ws.Rows(318).Insert()

# Populate the newly inserted row 318 with the new price-report record.
$ws.Range("A318").Value = 5
$ws.Range("B318").Value = "Macroferia Regional de Talca"
$ws.Range("C318").Value = "Maule"
$ws.Range("D318").Value = 44610
$ws.Range("D318").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E318").Value = 7
$ws.Range("F318").Value = 100112043
$ws.Range("G318").Value = "Pepino ensalada"
$ws.Range("H318").Value = "Sin especificar"
$ws.Range("I318").Value = "Primera"
$ws.Range("J318").Value = 300
$ws.Range("K318").Value = 13000
$ws.Range("L318").Value = 13000
$ws.Range("M318").Value = 13000
$ws.Range("N318").Value = "$/caja 80 unidades"
$ws.Range("O318").Value = "Región del Maule"
$ws.Range("P318").Value = 162
$ws.Range("Q318").Value = 80
$ws.Range("R318").Value = "Hortaliza"
